# "Generate Report for Handoff"
#
# Refreshes the handoff-report workbook with a newly generated package:
# the old GUID-named markdown/xliff artifacts are replaced with a freshly
# generated set (new GUID + new content hash), and the recorded handoff /
# handback timestamps move forward by a few seconds to reflect the new
# generation run.

$wb = $excel.ActiveWorkbook

$oldGuid = "6df3b3c5-4243-409f-809c-1616146669c7"
$newGuid = "27c9f8cf-b105-4f41-9e4b-7b5afba4372f"

$oldHash = "8b3ab849636dc5e5eb82d2c21f4f28c8a9490f0e"
$newHash = "3a8af80ed43cfd1a7573b360c511e20d9b1de1f8"

$oldMdName = $oldGuid + ".md"
$newMdName = $newGuid + ".md"

$oldMdPath = "e2e\" + $oldMdName
$newMdPath = "e2e\" + $newMdName

$zhCnOld = $oldGuid + "." + $oldHash + ".zh-cn.xlf"
$zhCnNew = $newGuid + "." + $newHash + ".zh-cn.xlf"

$deDeOld = $oldGuid + "." + $oldHash + ".de-de.xlf"
$deDeNew = $newGuid + "." + $newHash + ".de-de.xlf"

$hoDateOld = "2016-08-20 09:03:32"
$hoDateNew = "2016-08-20 09:03:47"

$handoffDateOld = "2016-08-20 09:03:28"
$handoffDateNew = "2016-08-20 09:03:44"

# The hyperlink targets (the GitHub blob URLs) keep pointing at the same
# commit/path that was already on record - only the displayed file name
# (and, separately, the cell text for non-linked cells) moves to the new
# generated file.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b55f58c6dfa93217a5eecd9af7f714b8ff4e40a5/e2e/" + $oldMdName

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: plain file name (no hyperlink on this cell)
$wsOverview.Range("A2").Value = $newMdName

# B2: path + file name, carries the hyperlink - preserve the existing
# hyperlink target, only change the text that is shown.
$linkB2 = $wsOverview.Range("B2")
$linkB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($linkB2, $hyperlinkAddress, "", "", $newMdPath)

# G2: latest HO xliff generate date
$wsOverview.Range("G2").Value = $hoDateNew

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2: source file name, carries the hyperlink
$linkZhA2 = $wsZhCn.Range("A2")
$linkZhA2.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($linkZhA2, $hyperlinkAddress, "", "", $newMdName)

# G2: latest handoff file
$wsZhCn.Range("G2").Value = $zhCnNew

# H2: latest handoff datetime
$wsZhCn.Range("H2").Value = $handoffDateNew

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2: source file name, carries the hyperlink
$linkDeA2 = $wsDeDe.Range("A2")
$linkDeA2.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($linkDeA2, $hyperlinkAddress, "", "", $newMdName)

# G2: latest handback file
$wsDeDe.Range("G2").Value = $deDeNew

# H2: latest handback datetime
$wsDeDe.Range("H2").Value = $hoDateNew
